# Applies scheduled-runner market-price/profit updates to each profession sheet's table.
# Generated from the canonical OOXML diff: updates currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H:N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 43.8
$ws.Range("I29").Value = 43.8
$ws.Range("K29").Value = 131.4
$ws.Range("M29").Value = 149.6
$ws.Range("H86").Value = 2062.8333
$ws.Range("I86").Value = 2083.3333
$ws.Range("K86").Value = 2083.3333
$ws.Range("M86").Value = -960.3332999999998
$ws.Range("H89").Value = 2062.8333
$ws.Range("I89").Value = 2083.3333
$ws.Range("K89").Value = 10416.6665
$ws.Range("M89").Value = -4800.666499999999
$ws.Range("H123").Value = 41875
$ws.Range("J123").Value = 41875
$ws.Range("L123").Value = 41875
$ws.Range("N123").Value = -51675
$ws.Range("H132").Value = 1904.6346
$ws.Range("I132").Value = 1685.25
$ws.Range("J132").Value = 3111.25
$ws.Range("K132").Value = 5055.75
$ws.Range("L132").Value = 9333.75
$ws.Range("M132").Value = -2525.75
$ws.Range("N132").Value = -14393.75
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 125300
$ws.Range("J134").Value = 125300
$ws.Range("L134").Value = 125300
$ws.Range("N134").Value = -135440
$ws.Range("H135").Value = 972.36
$ws.Range("I135").Value = 972.36
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8751.24
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6216.24
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1657.5555
$ws.Range("I137").Value = 1355.8572
$ws.Range("J137").Value = 2713.5
$ws.Range("K137").Value = 4067.5716
$ws.Range("L137").Value = 8140.5
$ws.Range("M137").Value = -1517.5716
$ws.Range("N137").Value = -13240.5
$ws.Range("H138").Value = 2919.6104
$ws.Range("I138").Value = 1923.1136
$ws.Range("J138").Value = 4248.273
$ws.Range("K138").Value = 5769.3408
$ws.Range("L138").Value = 12744.819
$ws.Range("M138").Value = -629.3407999999999
$ws.Range("N138").Value = -23024.819
$ws.Range("H140").Value = 86216.664
$ws.Range("J140").Value = 95960
$ws.Range("L140").Value = 95960
$ws.Range("N140").Value = -106320
$ws.Range("H141").Value = 3525.2222
$ws.Range("I141").Value = 1612.8823
$ws.Range("J141").Value = 36035
$ws.Range("K141").Value = 4838.6469
$ws.Range("L141").Value = 108105
$ws.Range("M141").Value = 341.3531000000003
$ws.Range("N141").Value = -118465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9026.393
$ws.Range("I32").Value = 9285.823
$ws.Range("K32").Value = 9285.823
$ws.Range("M32").Value = -8998.823
$ws.Range("H61").Value = 1328.7241
$ws.Range("I61").Value = 1343.0869
$ws.Range("J61").Value = 1273.6666
$ws.Range("K61").Value = 1343.0869
$ws.Range("L61").Value = 1273.6666
$ws.Range("M61").Value = -1131.0869
$ws.Range("N61").Value = -1697.6666
$ws.Range("H63").Value = 3353.6875
$ws.Range("I63").Value = 3435.3076
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 3435.3076
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -2749.3076
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 3353.6875
$ws.Range("I66").Value = 3435.3076
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 17176.538
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -13744.538
$ws.Range("N66").Value = -21864
$ws.Range("H74").Value = 746.58
$ws.Range("I74").Value = 729.2222
$ws.Range("J74").Value = 902.8
$ws.Range("K74").Value = 729.2222
$ws.Range("L74").Value = 902.8
$ws.Range("M74").Value = 144.7778
$ws.Range("N74").Value = -2650.8
$ws.Range("H77").Value = 746.58
$ws.Range("I77").Value = 729.2222
$ws.Range("J77").Value = 902.8
$ws.Range("K77").Value = 3646.111
$ws.Range("L77").Value = 4514
$ws.Range("M77").Value = 721.8889999999997
$ws.Range("N77").Value = -13250
$ws.Range("H102").Value = 78270.69500000001
$ws.Range("I102").Value = 1456.2727
$ws.Range("J102").Value = 500750
$ws.Range("K102").Value = 1456.2727
$ws.Range("L102").Value = 500750
$ws.Range("M102").Value = 165.7273
$ws.Range("N102").Value = -503994
$ws.Range("H119").Value = 37998.5
$ws.Range("J119").Value = 37998.5
$ws.Range("L119").Value = 37998.5
$ws.Range("N119").Value = -47674.5
$ws.Range("H125").Value = 67500.75
$ws.Range("J125").Value = 67500.75
$ws.Range("L125").Value = 67500.75
$ws.Range("N125").Value = -77340.75
$ws.Range("H136").Value = 1328.7241
$ws.Range("I136").Value = 1343.0869
$ws.Range("J136").Value = 1273.6666
$ws.Range("K136").Value = 4029.2607
$ws.Range("L136").Value = 3820.9998
$ws.Range("M136").Value = -1479.2607
$ws.Range("N136").Value = -8920.9998
$ws.Range("H139").Value = 59702.35
$ws.Range("J139").Value = 59496.25
$ws.Range("L139").Value = 59496.25
$ws.Range("N139").Value = -69776.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 18110
$ws.Range("J81").Value = 18110
$ws.Range("L81").Value = 18110
$ws.Range("N81").Value = -20232
$ws.Range("H84").Value = 18110
$ws.Range("J84").Value = 18110
$ws.Range("L84").Value = 54330
$ws.Range("N84").Value = -64938
$ws.Range("H134").Value = 2142.6304
$ws.Range("I134").Value = 1657.4193
$ws.Range("J134").Value = 3145.4
$ws.Range("K134").Value = 4972.257900000001
$ws.Range("L134").Value = 9436.200000000001
$ws.Range("M134").Value = -2437.257900000001
$ws.Range("N134").Value = -14506.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1602.6111
$ws.Range("I31").Value = 1145.3829
$ws.Range("K31").Value = 1145.3829
$ws.Range("M31").Value = -850.3829000000001
$ws.Range("H34").Value = 1602.6111
$ws.Range("I34").Value = 1145.3829
$ws.Range("K34").Value = 1145.3829
$ws.Range("M34").Value = -943.3829000000001
$ws.Range("H93").Value = 9644.583000000001
$ws.Range("I93").Value = 7483.9
$ws.Range("J93").Value = 20448
$ws.Range("K93").Value = 7483.9
$ws.Range("L93").Value = 20448
$ws.Range("M93").Value = -5611.9
$ws.Range("N93").Value = -24192
$ws.Range("H103").Value = 12104.5
$ws.Range("I103").Value = 6472.6665
$ws.Range("J103").Value = 29000
$ws.Range("K103").Value = 6472.6665
$ws.Range("L103").Value = 29000
$ws.Range("M103").Value = -5300.6665
$ws.Range("N103").Value = -31344
$ws.Range("H132").Value = 376775.6
$ws.Range("I132").Value = 423231.03
$ws.Range("J132").Value = 5132
$ws.Range("K132").Value = 1269693.09
$ws.Range("L132").Value = 15396
$ws.Range("M132").Value = -1267163.09
$ws.Range("N132").Value = -20456
$ws.Range("H134").Value = 1208.4
$ws.Range("I134").Value = 999.97437
$ws.Range("J134").Value = 1716.4375
$ws.Range("K134").Value = 2999.92311
$ws.Range("L134").Value = 5149.3125
$ws.Range("M134").Value = -464.9231100000002
$ws.Range("N134").Value = -10219.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2300
$ws.Range("I123").Value = 950
$ws.Range("J123").Value = 5000
$ws.Range("K123").Value = 2850
$ws.Range("L123").Value = 15000
$ws.Range("M123").Value = -400
$ws.Range("N123").Value = -19900
$ws.Range("H133").Value = 4639.4707
$ws.Range("I133").Value = 2586.2
$ws.Range("J133").Value = 5495
$ws.Range("K133").Value = 7758.599999999999
$ws.Range("L133").Value = 16485
$ws.Range("M133").Value = -2698.599999999999
$ws.Range("N133").Value = -26605
$ws.Range("H134").Value = 3459.55
$ws.Range("I134").Value = 1181.381
$ws.Range("J134").Value = 5977.5264
$ws.Range("K134").Value = 3544.143
$ws.Range("L134").Value = 17932.5792
$ws.Range("M134").Value = 1525.857
$ws.Range("N134").Value = -28072.5792
$ws.Range("H137").Value = 25643878
$ws.Range("J137").Value = 41669936
$ws.Range("L137").Value = 125009808
$ws.Range("N137").Value = -125020008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5750
$ws.Range("J43").Value = 5750
$ws.Range("L43").Value = 5750
$ws.Range("N43").Value = -6052
$ws.Range("H46").Value = 16450
$ws.Range("J46").Value = 19933.334
$ws.Range("L46").Value = 19933.334
$ws.Range("N46").Value = -20245.334
$ws.Range("H126").Value = 2042
$ws.Range("I126").Value = 1458.4
$ws.Range("J126").Value = 2307.2727
$ws.Range("K126").Value = 4375.200000000001
$ws.Range("L126").Value = 6921.8181
$ws.Range("M126").Value = -1905.200000000001
$ws.Range("N126").Value = -11861.8181
$ws.Range("H132").Value = 1533.6216
$ws.Range("I132").Value = 1204.5151
$ws.Range("J132").Value = 4248.75
$ws.Range("K132").Value = 3613.5453
$ws.Range("L132").Value = 12746.25
$ws.Range("M132").Value = -1083.5453
$ws.Range("N132").Value = -17806.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1534.5
$ws.Range("I68").Value = 1332.25
$ws.Range("J68").Value = 1837.875
$ws.Range("K68").Value = 1332.25
$ws.Range("L68").Value = 1837.875
$ws.Range("M68").Value = -583.25
$ws.Range("N68").Value = -3335.875
$ws.Range("H71").Value = 1534.5
$ws.Range("I71").Value = 1332.25
$ws.Range("J71").Value = 1837.875
$ws.Range("K71").Value = 6661.25
$ws.Range("L71").Value = 9189.375
$ws.Range("M71").Value = -2917.25
$ws.Range("N71").Value = -16677.375
$ws.Range("H93").Value = 945.6
$ws.Range("J93").Value = 1400
$ws.Range("L93").Value = 1400
$ws.Range("N93").Value = -3896
$ws.Range("H110").Value = 22620
$ws.Range("J110").Value = 22620
$ws.Range("L110").Value = 22620
$ws.Range("N110").Value = -30800
$ws.Range("H136").Value = 2117.0447
$ws.Range("I136").Value = 1821.6923
$ws.Range("J136").Value = 3140.9333
$ws.Range("K136").Value = 5465.0769
$ws.Range("L136").Value = 9422.7999
$ws.Range("M136").Value = -2915.0769
$ws.Range("N136").Value = -14522.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6722.2856
$ws.Range("I126").Value = 8021.9414
$ws.Range("J126").Value = 1198.75
$ws.Range("K126").Value = 24065.8242
$ws.Range("L126").Value = 3596.25
$ws.Range("M126").Value = -21595.8242
$ws.Range("N126").Value = -8536.25
$ws.Range("H132").Value = 792.5323
$ws.Range("I132").Value = 661.3137
$ws.Range("J132").Value = 1400.909
$ws.Range("K132").Value = 1983.9411
$ws.Range("L132").Value = 4202.727000000001
$ws.Range("M132").Value = 546.0589
$ws.Range("N132").Value = -9262.727000000001

